$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27772.361
$ws.Range("I32").Value = 6358.45
$ws.Range("J32").Value = 150137.58
$ws.Range("K32").Value = 6358.45
$ws.Range("L32").Value = 150137.58
$ws.Range("M32").Value = -6071.45
$ws.Range("N32").Value = -150711.58

$ws.Range("H44").Value = 25479.6
$ws.Range("J44").Value = 25479.6
$ws.Range("L44").Value = 25479.6
$ws.Range("N44").Value = -26455.6

$ws.Range("H55").Value = 21399
$ws.Range("J55").Value = 24624.166
$ws.Range("L55").Value = 24624.166
$ws.Range("N55").Value = -25254.166

$ws.Range("H61").Value = 2225.2
$ws.Range("I61").Value = 1884.9231
$ws.Range("K61").Value = 1884.9231
$ws.Range("M61").Value = -1672.9231

$ws.Range("H80").Value = 25947.4
$ws.Range("J80").Value = 25947.4
$ws.Range("L80").Value = 25947.4
$ws.Range("N80").Value = -27943.4

$ws.Range("H83").Value = 25947.4
$ws.Range("J83").Value = 25947.4
$ws.Range("L83").Value = 77842.20000000001
$ws.Range("N83").Value = -87826.20000000001

$ws.Range("H136").Value = 2225.2
$ws.Range("I136").Value = 1884.9231
$ws.Range("K136").Value = 5654.7693
$ws.Range("M136").Value = -3104.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34737.6
$ws.Range("J35").Value = 34737.6
$ws.Range("L35").Value = 34737.6
$ws.Range("N35").Value = -35357.6

$ws.Range("H94").Value = 723.7406999999999
$ws.Range("I94").Value = 598.55
$ws.Range("K94").Value = 598.55
$ws.Range("M94").Value = -147.55

$ws.Range("H105").Value = 2620.3
$ws.Range("I105").Value = 2133.6667
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 2133.6667
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -386.6667000000002
$ws.Range("N105").Value = -10494

$ws.Range("H107").Value = 7132.095
$ws.Range("I107").Value = 7462.4116
$ws.Range("J107").Value = 5728.25
$ws.Range("K107").Value = 7462.4116
$ws.Range("L107").Value = 5728.25
$ws.Range("M107").Value = -5542.4116
$ws.Range("N107").Value = -9568.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1524.6666
$ws.Range("I16").Value = 1444.4
$ws.Range("K16").Value = 1444.4
$ws.Range("M16").Value = -1157.4

$ws.Range("H60").Value = 13738.875
$ws.Range("J60").Value = 13738.875
$ws.Range("L60").Value = 13738.875
$ws.Range("N60").Value = -14760.875

$ws.Range("H86").Value = 83343160
$ws.Range("I86").Value = 111123120
$ws.Range("J86").Value = 3266.6667
$ws.Range("K86").Value = 111123120
$ws.Range("L86").Value = 3266.6667
$ws.Range("M86").Value = -111121997
$ws.Range("N86").Value = -5512.6667

$ws.Range("H89").Value = 83343160
$ws.Range("I89").Value = 111123120
$ws.Range("J89").Value = 3266.6667
$ws.Range("K89").Value = 555615600
$ws.Range("L89").Value = 16333.3335
$ws.Range("M89").Value = -555609984
$ws.Range("N89").Value = -27565.3335

$ws.Range("H105").Value = 4451.357
$ws.Range("I105").Value = 4265.364
$ws.Range("K105").Value = 4265.364
$ws.Range("M105").Value = -2518.364

$ws.Range("H113").Value = 1524.6666
$ws.Range("I113").Value = 1444.4
$ws.Range("K113").Value = 1444.4
$ws.Range("M113").Value = 725.5999999999999

$ws.Range("H134").Value = 7196.143
$ws.Range("I134").Value = 7305.95
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 21917.85
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -19382.85
$ws.Range("N134").Value = -20070

$ws.Range("H135").Value = 58347
$ws.Range("J135").Value = 58347
$ws.Range("L135").Value = 58347
$ws.Range("N135").Value = -68487

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 899.2857
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 891.53845
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 2674.61535
$ws.Range("M5").Value = -2888
$ws.Range("N5").Value = -2898.61535

$ws.Range("H122").Value = 2499.74
$ws.Range("J122").Value = 2997.805
$ws.Range("L122").Value = 26980.245
$ws.Range("N122").Value = -31880.245

$ws.Range("H134").Value = 4428.5713
$ws.Range("I134").Value = 4450
$ws.Range("J134").Value = 4400
$ws.Range("K134").Value = 13350
$ws.Range("L134").Value = 13200
$ws.Range("M134").Value = -8280
$ws.Range("N134").Value = -23340

$ws.Range("H135").Value = 899.2857
$ws.Range("I135").Value = 1000
$ws.Range("J135").Value = 891.53845
$ws.Range("K135").Value = 9000
$ws.Range("L135").Value = 8023.84605
$ws.Range("M135").Value = -6465
$ws.Range("N135").Value = -13093.84605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1863.5454
$ws.Range("I97").Value = 1874.875
$ws.Range("K97").Value = 1874.875
$ws.Range("M97").Value = -1378.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1729.4
$ws.Range("I61").Value = 1461.9584
$ws.Range("K61").Value = 1461.9584
$ws.Range("M61").Value = -1259.9584

$ws.Range("H82").Value = 1882.6875
$ws.Range("I82").Value = 1365
$ws.Range("J82").Value = 2400.375
$ws.Range("K82").Value = 1365
$ws.Range("L82").Value = 2400.375
$ws.Range("M82").Value = -1004
$ws.Range("N82").Value = -3122.375

$ws.Range("H85").Value = 1882.6875
$ws.Range("I85").Value = 1365
$ws.Range("J85").Value = 2400.375
$ws.Range("K85").Value = 1365
$ws.Range("L85").Value = 2400.375
$ws.Range("M85").Value = -117
$ws.Range("N85").Value = -4896.375

$ws.Range("H113").Value = 1729.4
$ws.Range("I113").Value = 1461.9584
$ws.Range("K113").Value = 1461.9584
$ws.Range("M113").Value = 708.0416

$ws.Range("H133").Value = 46464.75
$ws.Range("J133").Value = 46464.75
$ws.Range("L133").Value = 46464.75
$ws.Range("N133").Value = -51524.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4007.5
$ws.Range("I62").Value = 4350
$ws.Range("J62").Value = 3939
$ws.Range("K62").Value = 4350
$ws.Range("L62").Value = 3939
$ws.Range("M62").Value = -3726
$ws.Range("N62").Value = -5187

$ws.Range("H65").Value = 4007.5
$ws.Range("I65").Value = 4350
$ws.Range("J65").Value = 3939
$ws.Range("K65").Value = 21750
$ws.Range("L65").Value = 19695
$ws.Range("M65").Value = -18630
$ws.Range("N65").Value = -25935

$ws.Range("H109").Value = 24294.25
$ws.Range("J109").Value = 24294.25
$ws.Range("L109").Value = 24294.25
$ws.Range("N109").Value = -27068.25

$ws.Range("H136").Value = 1154.7037
$ws.Range("I136").Value = 1058.878
$ws.Range("K136").Value = 3176.634
$ws.Range("M136").Value = -626.634
